$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D5").Value = 26
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = 12.567
$ws.Range("G5").Value = 571
$ws.Range("H5").Value = 280
$ws.Range("I5").Value = 53
$ws.Range("J5").Value = 59
$ws.Range("K5").Value = 3246
$ws.Range("L5").Value = 67
$ws.Range("M5").Value = 124
$ws.Range("N5").Value = 208
$ws.Range("O5").Value = 2.576923076923077
$ws.Range("P5").Value = 4.769230769230769
$ws.Range("Q5").Value = 8
$ws.Range("R5").Value = 21.96
$ws.Range("S5").Value = 29

$ws.Range("D6").Value = 109
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = 54.66733333333335
$ws.Range("G6").Value = 2273
$ws.Range("H6").Value = 953
$ws.Range("I6").Value = 312
$ws.Range("J6").Value = 304
$ws.Range("K6").Value = 22466
$ws.Range("L6").Value = 944
$ws.Range("M6").Value = 682
$ws.Range("N6").Value = 842
$ws.Range("O6").Value = 8.660550458715596
$ws.Range("P6").Value = 6.256880733944954
$ws.Range("Q6").Value = 7.724770642201835
$ws.Range("R6").Value = 20.85
$ws.Range("S6").Value = 30.09

$ws.Range("D7").Value = 40
$ws.Range("E7").Value = 34
$ws.Range("F7").Value = 21.07966666666666
$ws.Range("G7").Value = 2044
$ws.Range("H7").Value = 971
$ws.Range("I7").Value = 274
$ws.Range("J7").Value = 266
$ws.Range("K7").Value = 4224
$ws.Range("L7").Value = 219
$ws.Range("M7").Value = 153
$ws.Range("N7").Value = 467
$ws.Range("O7").Value = 5.475
$ws.Range("P7").Value = 3.825
$ws.Range("Q7").Value = 11.675
$ws.Range("R7").Value = 51.1
$ws.Range("S7").Value = 31.62

$ws.Range("D10").Value = 21
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 10.70233333333334
$ws.Range("G10").Value = 1537
$ws.Range("H10").Value = 625
$ws.Range("I10").Value = 217
$ws.Range("J10").Value = 159
$ws.Range("K10").Value = 848
$ws.Range("L10").Value = 65
$ws.Range("M10").Value = 153
$ws.Range("N10").Value = 273
$ws.Range("O10").Value = 3.095238095238095
$ws.Range("P10").Value = 7.285714285714286
$ws.Range("Q10").Value = 13
$ws.Range("R10").Value = 73.19
$ws.Range("S10").Value = 30.58

$ws.Range("D12").Value = 61
$ws.Range("F12").Value = 30.73633333333333
$ws.Range("G12").Value = 4895
$ws.Range("H12").Value = 2793
$ws.Range("I12").Value = 609
$ws.Range("J12").Value = 662
$ws.Range("K12").Value = 1175
$ws.Range("L12").Value = 148
$ws.Range("M12").Value = 181
$ws.Range("N12").Value = 1013
$ws.Range("O12").Value = 2.426229508196721
$ws.Range("P12").Value = 2.967213114754098
$ws.Range("Q12").Value = 16.60655737704918
$ws.Range("R12").Value = 80.25
$ws.Range("S12").Value = 30.23

$ws.Range("D13").Value = 38
$ws.Range("E13").Value = 23
$ws.Range("F13").Value = 19.95216666666667
$ws.Range("G13").Value = 1209
$ws.Range("H13").Value = 512
$ws.Range("I13").Value = 105
$ws.Range("J13").Value = 197
$ws.Range("K13").Value = 7475
$ws.Range("L13").Value = 242
$ws.Range("M13").Value = 198
$ws.Range("N13").Value = 289
$ws.Range("O13").Value = 6.368421052631579
$ws.Range("P13").Value = 5.210526315789473
$ws.Range("Q13").Value = 7.605263157894737
$ws.Range("R13").Value = 31.82
$ws.Range("S13").Value = 31.5
